$wb = $excel.ActiveWorkbook

# --- charge-order sheet: move selection from C30 to B17, deactivate tab ---
$wsCharge = $wb.Worksheets.Item("charge-order")
$wsCharge.Range("B17").Select()

# --- cash-in cash-out sheet: remove the "password" field row (row 48) from the
#     "Create order topup" API section; this shifts all subsequent rows up by one ---
$wsCashInOut = $wb.Worksheets.Item("cash-in cash-out")
$wsCashInOut.Rows.Item(48).Delete()

# make this sheet the active tab and move the selection to C14
$wsCashInOut.Activate()
$wsCashInOut.Range("C14").Select()
